# Home Page Content.docx edit script
# 1) "Here, I am devotedly studying of computer game technology, ..." ->
#    "Here, I am" | " " | "studying computer game technology, ..."
# 2) "To gain a deeper understanding ... comprehensive details." ->
#    "For further details of my" | " education and professional experience, " |
#    "please refer to the CV provided below." and the paragraph break that used
#    to follow is converted into a line break (merging the "About" paragraph
#    with the "Skills" paragraph).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: fix wording around "Here, I am ... studying ..."
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute(
    "Here, I am devotedly studying of computer game technology, to expand my expertise and deepen my understanding.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Here, I am studying computer game technology, to expand my expertise and deepen my understanding.",
    2)

# Split "Here, I am" | " " | "studying computer game technology, ..." into
# three separate runs (matching the target markup) by forcing a formatting
# round-trip on the relevant sub-ranges.
$rng = $d.Content
$null = $rng.Find.Execute("Here, I am", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0

$rng = $d.Content
$null = $rng.Find.Execute(" studying computer game technology, to expand my expertise and deepen my understanding.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spaceRng = $d.Range($rng.Start, $rng.Start + 1)
$spaceRng.Bold = 1
$spaceRng.Bold = 0

# ---------------------------------------------------------------------------
# Edit 2: reword the CV sentence
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute(
    "To gain a deeper understanding of my education and professional experience, I encourage you to review the attached CV for comprehensive details.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "For further details of my education and professional experience, please refer to the CV provided below.",
    2)

# Split into three runs: "For further details of my" | " education and
# professional experience, " | "please refer to the CV provided below."
$rng = $d.Content
$null = $rng.Find.Execute("For further details of my", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0

$rng = $d.Content
$null = $rng.Find.Execute("please refer to the CV provided below.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1
$rng.Bold = 0

# ---------------------------------------------------------------------------
# Merge the "About" paragraph with the following "Skills" paragraph: the
# trailing paragraph mark becomes a simple line break (the "Skills" paragraph
# already starts with its own <w:br/>, so we end up with two consecutive
# breaks, exactly as in the target markup).
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("please refer to the CV provided below.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markPos = $rng.End
$markRng = $d.Range($markPos, $markPos + 1)
$markRng.Delete()
$insRng = $d.Range($markPos, $markPos)
$insRng.InsertBreak(6)  # wdLineBreak

Write-Output $d.Paragraphs.Item(2).Range.Text
